$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4833.3
$ws.Range("I62").Value = 4763.125
$ws.Range("K62").Value = 4763.125
$ws.Range("M62").Value = -4139.125
$ws.Range("H65").Value = 4833.3
$ws.Range("I65").Value = 4763.125
$ws.Range("K65").Value = 23815.625
$ws.Range("M65").Value = -20695.625
$ws.Range("H70").Value = 1496.1
$ws.Range("I70").Value = 810
$ws.Range("J70").Value = 1953.5
$ws.Range("K70").Value = 2430
$ws.Range("L70").Value = 5860.5
$ws.Range("M70").Value = -2160
$ws.Range("N70").Value = -6400.5
$ws.Range("H73").Value = 1496.1
$ws.Range("I73").Value = 810
$ws.Range("J73").Value = 1953.5
$ws.Range("K73").Value = 2430
$ws.Range("L73").Value = 5860.5
$ws.Range("M73").Value = -1494
$ws.Range("N73").Value = -7732.5
$ws.Range("H112").Value = 3176.8462
$ws.Range("J112").Value = 3176.8462
$ws.Range("L112").Value = 9530.5386
$ws.Range("N112").Value = -11746.5386
$ws.Range("H138").Value = 10460.077
$ws.Range("I138").Value = 17930.5
$ws.Range("J138").Value = 4056.8572
$ws.Range("K138").Value = 53791.5
$ws.Range("L138").Value = 12170.5716
$ws.Range("M138").Value = -48651.5
$ws.Range("N138").Value = -22450.5716
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2450.9058
$ws.Range("I32").Value = 1027.6438
$ws.Range("J32").Value = 11109.083
$ws.Range("K32").Value = 1027.6438
$ws.Range("L32").Value = 11109.083
$ws.Range("M32").Value = -740.6438000000001
$ws.Range("N32").Value = -11683.083
$ws.Range("I61").Value = 30717.361
$ws.Range("K61").Value = 30717.361
$ws.Range("M61").Value = -30505.361
$ws.Range("H74").Value = 430762.75
$ws.Range("I74").Value = 1765
$ws.Range("K74").Value = 1765
$ws.Range("M74").Value = -891
$ws.Range("H77").Value = 430762.75
$ws.Range("I77").Value = 1765
$ws.Range("K77").Value = 8825
$ws.Range("M77").Value = -4457
$ws.Range("H88").Value = 1059
$ws.Range("J88").Value = 1054.4
$ws.Range("L88").Value = 1054.4
$ws.Range("N88").Value = -1866.4
$ws.Range("H91").Value = 1059
$ws.Range("J91").Value = 1054.4
$ws.Range("L91").Value = 1054.4
$ws.Range("N91").Value = -3862.4
$ws.Range("H132").Value = 1938.4667
$ws.Range("I132").Value = 1581.1945
$ws.Range("J132").Value = 3367.5557
$ws.Range("K132").Value = 4743.583500000001
$ws.Range("L132").Value = 10102.6671
$ws.Range("M132").Value = -2213.583500000001
$ws.Range("N132").Value = -15162.6671
$ws.Range("I136").Value = 30717.361
$ws.Range("K136").Value = 92152.083
$ws.Range("M136").Value = -89602.083
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8389.543
$ws.Range("I99").Value = 7185.2964
$ws.Range("K99").Value = 7185.2964
$ws.Range("M99").Value = -5687.2964
$ws.Range("H134").Value = 17309814
$ws.Range("I134").Value = 1861.2727
$ws.Range("K134").Value = 5583.8181
$ws.Range("M134").Value = -3048.8181
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1758.9841
$ws.Range("I31").Value = 1152.0312
$ws.Range("K31").Value = 1152.0312
$ws.Range("M31").Value = -857.0311999999999
$ws.Range("H34").Value = 1758.9841
$ws.Range("I34").Value = 1152.0312
$ws.Range("K34").Value = 1152.0312
$ws.Range("M34").Value = -950.0311999999999
$ws.Range("H105").Value = 1344.6666
$ws.Range("I105").Value = 1519.9231
$ws.Range("K105").Value = 1519.9231
$ws.Range("M105").Value = 227.0769
$ws.Range("H132").Value = 19609816
$ws.Range("I132").Value = 1932.6666
$ws.Range("K132").Value = 5797.9998
$ws.Range("M132").Value = -3267.9998
$ws.Range("H134").Value = 1713.7878
$ws.Range("I134").Value = 1556.7084
$ws.Range("J134").Value = 2132.6667
$ws.Range("K134").Value = 4670.1252
$ws.Range("L134").Value = 6398.000100000001
$ws.Range("M134").Value = -2135.1252
$ws.Range("N134").Value = -11468.0001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 30480200
$ws.Range("I129").Value = 6522111
$ws.Range("K129").Value = 19566333
$ws.Range("M129").Value = -19561333
$ws.Range("H131").Value = 6124709
$ws.Range("I131").Value = 7577174
$ws.Range("K131").Value = 22731522
$ws.Range("M131").Value = -22726482
$ws.Range("H134").Value = 1953.1364
$ws.Range("I134").Value = 1427.0952
$ws.Range("K134").Value = 4281.2856
$ws.Range("M134").Value = 788.7143999999998
$ws.Range("H139").Value = 4562.143
$ws.Range("I139").Value = 4221.615
$ws.Range("K139").Value = 12664.845
$ws.Range("M139").Value = -7524.844999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5766
$ws.Range("I70").Value = 5740.75
$ws.Range("J70").Value = 5799.6665
$ws.Range("K70").Value = 5740.75
$ws.Range("L70").Value = 5799.6665
$ws.Range("M70").Value = -5470.75
$ws.Range("N70").Value = -6339.6665
$ws.Range("H73").Value = 5766
$ws.Range("I73").Value = 5740.75
$ws.Range("J73").Value = 5799.6665
$ws.Range("K73").Value = 5740.75
$ws.Range("L73").Value = 5799.6665
$ws.Range("M73").Value = -4804.75
$ws.Range("N73").Value = -7671.6665
$ws.Range("H80").Value = 23928.945
$ws.Range("I80").Value = 18374.77
$ws.Range("J80").Value = 38369.8
$ws.Range("K80").Value = 18374.77
$ws.Range("L80").Value = 38369.8
$ws.Range("M80").Value = -17376.77
$ws.Range("N80").Value = -40365.8
$ws.Range("H83").Value = 23928.945
$ws.Range("I83").Value = 18374.77
$ws.Range("J83").Value = 38369.8
$ws.Range("K83").Value = 91873.85000000001
$ws.Range("L83").Value = 191849
$ws.Range("M83").Value = -86881.85000000001
$ws.Range("N83").Value = -201833
$ws.Range("H97").Value = 883.2727
$ws.Range("J97").Value = 1499.5
$ws.Range("L97").Value = 1499.5
$ws.Range("N97").Value = -2491.5
$ws.Range("H102").Value = 11906691
$ws.Range("I102").Value = 12822418
$ws.Range("K102").Value = 12822418
$ws.Range("M102").Value = -12820796
$ws.Range("H126").Value = 2750.353
$ws.Range("I126").Value = 1722.4445
$ws.Range("K126").Value = 5167.333500000001
$ws.Range("M126").Value = -2697.333500000001
$ws.Range("H132").Value = 5629974.5
$ws.Range("I132").Value = 5097.1465
$ws.Range("K132").Value = 15291.4395
$ws.Range("M132").Value = -12761.4395
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8041.625
$ws.Range("I22").Value = 6537.375
$ws.Range("K22").Value = 6537.375
$ws.Range("M22").Value = -6242.375
$ws.Range("H27").Value = 8041.625
$ws.Range("I27").Value = 6537.375
$ws.Range("K27").Value = 6537.375
$ws.Range("M27").Value = -6430.375
$ws.Range("H55").Value = 1916.875
$ws.Range("I55").Value = 2446.25
$ws.Range("K55").Value = 2446.25
$ws.Range("M55").Value = -2273.25
$ws.Range("H61").Value = 3564.5833
$ws.Range("I61").Value = 2097.2222
$ws.Range("J61").Value = 7966.6665
$ws.Range("K61").Value = 2097.2222
$ws.Range("L61").Value = 7966.6665
$ws.Range("M61").Value = -1895.2222
$ws.Range("N61").Value = -8370.666499999999
$ws.Range("H113").Value = 3564.5833
$ws.Range("I113").Value = 2097.2222
$ws.Range("J113").Value = 7966.6665
$ws.Range("K113").Value = 2097.2222
$ws.Range("L113").Value = 7966.6665
$ws.Range("M113").Value = 72.77779999999984
$ws.Range("N113").Value = -12306.6665
$ws.Range("H122").Value = 3207.6924
$ws.Range("I122").Value = 2973.913
$ws.Range("K122").Value = 8921.739
$ws.Range("M122").Value = -6471.739
$ws.Range("H132").Value = 2834.6182
$ws.Range("I132").Value = 2484.9429
$ws.Range("K132").Value = 7454.8287
$ws.Range("M132").Value = -4924.8287
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 10677.777
$ws.Range("I14").Value = 19099.666
$ws.Range("K14").Value = 19099.666
$ws.Range("M14").Value = -18931.666
$ws.Range("H41").Value = 17079
$ws.Range("J41").Value = 17516.8
$ws.Range("L41").Value = 17516.8
$ws.Range("N41").Value = -18296.8
$ws.Range("H100").Value = 654738.8
$ws.Range("I100").Value = 785395.8
$ws.Range("J100").Value = 1453.6
$ws.Range("K100").Value = 1570791.6
$ws.Range("L100").Value = 2907.2
$ws.Range("M100").Value = -1570250.6
$ws.Range("N100").Value = -3989.2
$ws.Range("H122").Value = 3375.7334
$ws.Range("I122").Value = 2790.6365
$ws.Range("K122").Value = 8371.9095
$ws.Range("M122").Value = -5921.9095
$ws.Range("H132").Value = 1232.7941
$ws.Range("I132").Value = 923.8461
$ws.Range("J132").Value = 2236.875
$ws.Range("K132").Value = 2771.5383
$ws.Range("L132").Value = 6710.625
$ws.Range("M132").Value = -241.5383000000002
$ws.Range("N132").Value = -11770.625
$ws.Range("H136").Value = 5180.4644
$ws.Range("I136").Value = 6168.125
$ws.Range("K136").Value = 18504.375
$ws.Range("M136").Value = -15954.375
